# Generate Report for Handoff
# The file "a1bbab1d-6166-4e20-b695-03e46fdb686e" has just become ready for
# handoff, while "37391360-f6d2-46a2-9d14-f739163c26b9" was already ready.
# The status report re-sorts/refreshes rows 7 and 8 on every sheet so the
# two records trade places, and a1bbab1d's handoff timestamp is refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A7").Value = "37391360-f6d2-46a2-9d14-f739163c26b9.md"
$ws.Range("B7").Value = "Ready for handoff"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "2016-19-17 06:19:45"

$ws.Range("A8").Value = "a1bbab1d-6166-4e20-b695-03e46fdb686e.md"
$ws.Range("B8").Value = "Ready for handoff"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("D8").Value = "2016-23-17 06:23:15"

$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/2f1a226e470bb7cbfe67e51c104580dfe7015895/e2e/37391360-f6d2-46a2-9d14-f739163c26b9.md", "", "", "37391360-f6d2-46a2-9d14-f739163c26b9.md")
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/6f916d96a9ab0dfe75264bb693f7e7ec4154ad12/e2e/a1bbab1d-6166-4e20-b695-03e46fdb686e.md", "", "", "a1bbab1d-6166-4e20-b695-03e46fdb686e.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A7").Value = "37391360-f6d2-46a2-9d14-f739163c26b9.md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "37391360-f6d2-46a2-9d14-f739163c26b9.6a6c4fcf40381ced0545dad67ad536731dfc4ce0.zh-cn.xlf"
$ws.Range("E7").Value = "2016-03-17 06:19:37"

$ws.Range("A8").Value = "a1bbab1d-6166-4e20-b695-03e46fdb686e.md"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("D8").Value = "a1bbab1d-6166-4e20-b695-03e46fdb686e.2d3fb4c73e31de57fd1187938b579e26758e7a88.zh-cn.xlf"
$ws.Range("E8").Value = "2016-03-17 06:23:08"

$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/2f1a226e470bb7cbfe67e51c104580dfe7015895/e2e/37391360-f6d2-46a2-9d14-f739163c26b9.md", "", "", "37391360-f6d2-46a2-9d14-f739163c26b9.md")
$ws.Hyperlinks.Add($ws.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/09fccde9bc1c80948d3fce426706337f9f2ae155/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/37391360-f6d2-46a2-9d14-f739163c26b9.6a6c4fcf40381ced0545dad67ad536731dfc4ce0.zh-cn.xlf", "", "", "37391360-f6d2-46a2-9d14-f739163c26b9.6a6c4fcf40381ced0545dad67ad536731dfc4ce0.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/6f916d96a9ab0dfe75264bb693f7e7ec4154ad12/e2e/a1bbab1d-6166-4e20-b695-03e46fdb686e.md", "", "", "a1bbab1d-6166-4e20-b695-03e46fdb686e.md")
$ws.Hyperlinks.Add($ws.Range("D8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4a870f15162d94de5c6681eefa1bffd03213e058/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/a1bbab1d-6166-4e20-b695-03e46fdb686e.2d3fb4c73e31de57fd1187938b579e26758e7a88.zh-cn.xlf", "", "", "a1bbab1d-6166-4e20-b695-03e46fdb686e.2d3fb4c73e31de57fd1187938b579e26758e7a88.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A7").Value = "37391360-f6d2-46a2-9d14-f739163c26b9.md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "37391360-f6d2-46a2-9d14-f739163c26b9.6a6c4fcf40381ced0545dad67ad536731dfc4ce0.de-de.xlf"
$ws.Range("E7").Value = "2016-03-17 06:19:45"

$ws.Range("A8").Value = "a1bbab1d-6166-4e20-b695-03e46fdb686e.md"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("D8").Value = "a1bbab1d-6166-4e20-b695-03e46fdb686e.2d3fb4c73e31de57fd1187938b579e26758e7a88.de-de.xlf"
$ws.Range("E8").Value = "2016-03-17 06:23:15"

$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/2f1a226e470bb7cbfe67e51c104580dfe7015895/e2e/37391360-f6d2-46a2-9d14-f739163c26b9.md", "", "", "37391360-f6d2-46a2-9d14-f739163c26b9.md")
$ws.Hyperlinks.Add($ws.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/09fccde9bc1c80948d3fce426706337f9f2ae155/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/37391360-f6d2-46a2-9d14-f739163c26b9.6a6c4fcf40381ced0545dad67ad536731dfc4ce0.de-de.xlf", "", "", "37391360-f6d2-46a2-9d14-f739163c26b9.6a6c4fcf40381ced0545dad67ad536731dfc4ce0.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/6f916d96a9ab0dfe75264bb693f7e7ec4154ad12/e2e/a1bbab1d-6166-4e20-b695-03e46fdb686e.md", "", "", "a1bbab1d-6166-4e20-b695-03e46fdb686e.md")
$ws.Hyperlinks.Add($ws.Range("D8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4a870f15162d94de5c6681eefa1bffd03213e058/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/a1bbab1d-6166-4e20-b695-03e46fdb686e.2d3fb4c73e31de57fd1187938b579e26758e7a88.de-de.xlf", "", "", "a1bbab1d-6166-4e20-b695-03e46fdb686e.2d3fb4c73e31de57fd1187938b579e26758e7a88.de-de.xlf")
